# Updates the "cryptos" worksheet with refreshed price / volume data.
# Mirrors the automated "Updated cryptos list ... with GitHub Actions" commit:
# - refreshes Price (column D) and Volume(1h) (column E) figures for most rows
# - for rows 41/43 and 51, the coin identity (name/link) also changed along
#   with its price/volume (the feed re-ranked coins between runs)
#
# All values are written with a leading apostrophe so Excel stores them as
# literal text (matching the sheet's existing inline-string cell type)
# instead of auto-converting numeric-looking strings (e.g. "592.06") into
# numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '''66.690.11'
    'E2' = '''  +4.40%  '
    'D3' = '''3.492.35'
    'E3' = '''  +2.84%  '
    'D5' = '''592.06'
    'E5' = '''  +3.81%  '
    'D6' = '''169.57'
    'E6' = '''  +4.39%  '
    'E7' = '''  -0.10%  '
    'D8' = '''3.492.96'
    'E8' = '''  +2.82%  '
    'D9' = '''0.593'
    'E9' = '''  +8.51%  '
    'E10' = '''  +0.31%  '
    'E11' = '''  +7.36%  '
    'E12' = '''  +4.42%  '
    'D13' = '''4.088.76'
    'E13' = '''  +2.67%  '
    'E14' = '''  -0.19%  '
    'D15' = '''28.19'
    'E15' = '''  +4.78%  '
    'D16' = '''0.0000179'
    'E16' = '''  +4.26%  '
    'D17' = '''66.666.21'
    'E17' = '''  +4.26%  '
    'D18' = '''3.487.65'
    'E18' = '''  +2.59%  '
    'D19' = '''6.32'
    'E19' = '''  +3.62%  '
    'D20' = '''14.09'
    'E20' = '''  +4.10%  '
    'D21' = '''392.76'
    'E21' = '''  +4.38%  '
    'D22' = '''7.94'
    'E22' = '''  +2.26%  '
    'D23' = '''73.09'
    'E23' = '''  +4.24%  '
    'D25' = '''0.535'
    'E25' = '''  +4.94%  '
    'D26' = '''0.0000122'
    'E26' = '''  +7.44%  '
    'D27' = '''10.24'
    'E27' = '''  +7.66%  '
    'D28' = '''0.182'
    'E28' = '''  +1.70%  '
    'D29' = '''0.999'
    'E29' = '''  -0.08%  '
    'D30' = '''6.35'
    'E30' = '''  +4.57%  '
    'D31' = '''1.47'
    'E31' = '''  +5.91%  '
    'D32' = '''2.07'
    'E32' = '''  +3.70%  '
    'D33' = '''23.58'
    'E33' = '''  +3.64%  '
    'E34' = '''  +5.83%  '
    'E35' = '''  +0.14%  '
    'D36' = '''1.62'
    'E36' = '''  +9.92%  '
    'D37' = '''161.62'
    'E37' = '''  +1.29%  '
    'D38' = '''0.904'
    'E38' = '''  +4.95%  '
    'D39' = '''1.92'
    'E39' = '''  +6.77%  '
    'D40' = '''4.66'
    'E40' = '''  +7.01%  '
    'B41' = '''EnergySwap'
    'C41' = '''https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D41' = '''26.64'
    'E41' = '''  +3.67%  '
    'D42' = '''6.74'
    'E42' = '''  +5.29%  '
    'B43' = '''Hedera'
    'C43' = '''https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'D43' = '''0.0744'
    'E43' = '''  +3.27%  '
    'D44' = '''26.87'
    'E44' = '''  +2.68%  '
    'D45' = '''43.20'
    'E45' = '''  +1.12%  '
    'D46' = '''2.767.92'
    'E46' = '''  +1.21%  '
    'D47' = '''0.0313'
    'E47' = '''  +2.63%  '
    'D48' = '''2.49'
    'E48' = '''  +3.98%  '
    'D49' = '''346.31'
    'E49' = '''  +5.71%  '
    'E50' = '''  +5.65%  '
    'B51' = '''SuiNetwork'
    'C51' = '''https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
    'D51' = '''0.887'
    'E51' = '''  +9.04%  '
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
